$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '76.294.21'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +1.37%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.877.81'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +7.60%  '
$ws.Range('E4').Value = '  +0.12%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '196.24'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +4.64%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '599.77'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +2.14%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.00'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.14%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.553'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +3.66%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.193'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.22%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '2.872.43'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +7.41%  '
$ws.Range('E11').Value = '  +9.60%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.160'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -2.02%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.92'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +4.13%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '3.408.37'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +7.55%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '76.216.93'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.54%  '
$ws.Range('B16').Value = 'Avalanche'
$ws.Range('C16').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '27.51'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +3.95%  '
$ws.Range('B17').Value = 'ShibaInu'
$ws.Range('C17').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.0000189'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +1.43%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.880.55'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +7.17%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '9.04'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.30%  '
$ws.Range('E20').Value = '  +5.37%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '383.27'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +3.10%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '2.33'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +3.06%  '
$ws.Range('E23').Value = '  +1.36%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '71.75'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +2.71%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.999'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.02%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.026.53'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +8.28%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '4.22'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.97%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.74'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +4.08%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.0000105'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +10.77%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.00'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.09%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.40'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.48%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '513.25'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.27%  '
$ws.Range('E33').Value = '  +0.86%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.81'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +3.65%  '
$ws.Range('E35').Value = '  +0.14%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '167.13'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +2.72%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '20.05'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +4.37%  '
$ws.Range('E38').Value = '  +0.12%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '19.51'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.84%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '185.39'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +8.79%  '
$ws.Range('E41').Value = '  -0.08%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.345'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +4.44%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.07'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +1.44%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.67'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.55%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0915'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +8.10%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.23'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +3.85%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '40.27'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +3.07%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.37'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.49%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.578'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +9.01%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.681'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +15.29%  '
$ws.Range('E51').Value = '  +2.97%  '
